$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.722.95'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '3.485.57'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'593.14"
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").Value = "'171.45"
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = "'0.130"
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = "'0.431"
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '4.093.36'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = "'28.74"
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = '66.792.38'
$ws.Range("E15").Value = '  +0.23%  '
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '3.469.09'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = "'6.27"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = "'14.01"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").Value = "'391.97"
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = "'7.94"
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").Value = "'72.66"
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = "'0.533"
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").Value = "'5.69"
$ws.Range("E25").Value = '  -3.19%  '
$ws.Range("D26").Value = "'0.0000119"
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = "'6.20"
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = '  -3.55%  '
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").Value = "'23.62"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'7.30"
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("D39").Value = "'6.91"
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("D40").Value = "'4.65"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").Value = "'27.24"
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").Value = "'0.0739"
$ws.Range("E42").Value = '  -0.97%  '
$ws.Range("D43").Value = "'26.12"
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").Value = '2.800.00'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").Value = "'42.63"
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").Value = "'2.54"
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("D47").Value = "'0.0301"
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D48").Value = "'336.01"
$ws.Range("E48").Value = '  -4.43%  '
$ws.Range("D49").Value = "'34.24"
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("E50").Value = '  -3.30%  '
$ws.Range("D51").Value = "'0.104"
$ws.Range("E51").Value = '  -1.91%  '
